$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.947.60'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.36%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.234.22'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.16%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.52'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.625'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.88%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.77'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -6.39%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.407'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.18'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0902'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.67%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.74%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.565.35'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.55'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.78'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.67'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.75%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.67%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.239.88'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.847.73'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.59%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0912'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.70'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.61%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.20'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '248.61'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.11%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.11%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.32'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.74'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.31%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '169.71'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.142'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.35%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.24%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.49%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.55'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -10.37%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.56%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.92%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.91%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.67%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -8.37%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.90%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.63'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.000245'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.57%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.68'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.72%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.52'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -5.69%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '99.15'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0965'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.56%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.473.25'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.51%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.70'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -6.46%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.29'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +8.10%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.58%  '
